$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Terrains sheet: insert a new row 5 ("Walkable Map Edge" / id 2) and shift
# all the following terrain rows down by one.
# ---------------------------------------------------------------------------
$terrains = $wb.Worksheets.Item("Terrains")
$terrains.Rows.Item(5).Insert()
$terrains.Cells.Item(5, 2).Value = 2
$terrains.Cells.Item(5, 3).Value = "Walkable Map Edge"
$terrains.Cells.Item(5, 5).Value = 0
# Restore the formatting (bold, right aligned, quote-prefixed) that the row
# below (now row 7, style "19") uses, since writing the new values resets it.
$terrains.Cells.Item(7, 2).Copy()
$terrains.Cells.Item(5, 2).PasteSpecial(-4122)
$terrains.Range("E5").Select()

# ---------------------------------------------------------------------------
# Features sheet: just a cursor/selection move.
# ---------------------------------------------------------------------------
$features = $wb.Worksheets.Item("Features")
$features.Range("D16").Select()

# ---------------------------------------------------------------------------
# Units sheet: just a cursor/selection move.
# ---------------------------------------------------------------------------
$units = $wb.Worksheets.Item("Units")
$units.Range("J19").Select()

# ---------------------------------------------------------------------------
# Items sheet: clear E193 (the "2" inventory-capacity value is no longer
# present) and move the cursor/selection.
# ---------------------------------------------------------------------------
$items = $wb.Worksheets.Item("Items")
$items.Range("E193").ClearContents()
$items.Range("F191").Select()

# Items was (and remains) the active sheet/tab - make sure it stays active
# after touching the other sheets above.
$items.Activate()
